{"js": "// Apply strikethrough formatting to the four \"completed\" requirement\n// bullets (paragraph mark + every run), matching the author's edit that\n// crosses these items off as done while leaving the \"Optional\" section and\n// the \"unprofitable products\" bullet untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Identify target paragraphs by their (trimmed) text so the script is\n// resilient to any incidental paragraph re-indexing.\nconst targets = [\n  \"Product return quantities by geographic location (region, state)\",\n  \"Create a dashboard using two or more of the above data visualizations (\\u201cSheets\\u201d), with \\u201cUse as filter\\u201d turned on so that the visualizations are interactive (clicking on one, causes a reaction in the other)\",\n  \"Which products have the highest return rates?\",\n  \"Sometimes a product line needs to be discontinued. Other times, a product simply needs to be discontinued in certain geographies. Can you find a product that this applies to?\"\n];\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const text = (items[i].text || \"\").trim().replace(/\\s+/g, \" \");\n  if (targets.some((t) => text === t)) {\n    // Setting strikeThrough on the paragraph's Font applies it to the\n    // paragraph mark run properties AND every run of text in the paragraph,\n    // matching Word's \"select paragraph, toggle Strikethrough\" behavior.\n    items[i].font.strikeThrough = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to the four \"completed\" requirement\n# bullets (paragraph mark + every run), matching the author's edit that\n# crosses these items off as done while leaving the \"Optional\" section and\n# the \"unprofitable products\" bullet untouched.\n\n$d = $word.ActiveDocument\n\n$lq = [char]8220   # U+201C LEFT DOUBLE QUOTATION MARK\n$rq = [char]8221   # U+201D RIGHT DOUBLE QUOTATION MARK\n\n$targets = @(\n    \"Product return quantities by geographic location (region, state)\",\n    (\"Create a dashboard using two or more of the above data visualizations ({0}Sheets{1}), with {0}Use as filter{1} turned on so that the visualizations are interactive (clicking on one, causes a reaction in the other)\" -f $lq, $rq),\n    \"Which products have the highest return rates?\",\n    \"Sometimes a product line needs to be discontinued. Other times, a product simply needs to be discontinued in certain geographies. Can you find a product that this applies to?\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $norm = $p.Range.Text.Trim()\n    foreach ($t in $targets) {\n        if ($norm -eq $t) {\n            # Setting StrikeThrough on the paragraph Range's Font applies it to\n            # the paragraph mark run properties AND every run of text in the\n            # paragraph, matching Word's \"select paragraph, toggle Strikethrough\".\n            $p.Range.Font.StrikeThrough = 1\n        }\n    }\n}\n"}
